# Add OLE_LINK bookmarks around two "Opis testu" cell paragraphs.
$d = $word.ActiveDocument

# --- First location: "Poprawna i kompletna rejestracja użytkownika do Systemu
#     Zarządzania Fakturami" (TC_USER_REGISTER_001 - successful registration)
#     Wrap the whole paragraph text with OLE_LINK1 / OLE_LINK2.

$r1 = $d.Content
$r1.Find.Execute("Poprawna i kompletna rejestracja użytkownika do Systemu Zarządzania Fakturami", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$d.Bookmarks.Add("OLE_LINK1", $r1)

$r2 = $d.Content
$r2.Find.Execute("Poprawna i kompletna rejestracja użytkownika do Systemu Zarządzania Fakturami", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$d.Bookmarks.Add("OLE_LINK2", $r2)

# --- Second location: "Niepowodzenie operacji rejestracji użytkownika do
#     Systemu Zarządzania Fakturami - użytkownik już istnieje"
#     (TC_USER_REGISTER_00x - registration failure, user already exists)
#     Wrap the whole paragraph text with OLE_LINK3 / OLE_LINK4.

$r3 = $d.Content
$r3.Find.Execute("Niepowodzenie operacji rejestracji użytkownika do Systemu Zarządzania Fakturami - użytkownik już istnieje", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$d.Bookmarks.Add("OLE_LINK3", $r3)

$r4 = $d.Content
$r4.Find.Execute("Niepowodzenie operacji rejestracji użytkownika do Systemu Zarządzania Fakturami - użytkownik już istnieje", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$d.Bookmarks.Add("OLE_LINK4", $r4)
